$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 113778.11
$ws.Range("I76").Value = 113778.11
$ws.Range("K76").Value = 113778.11
$ws.Range("M76").Value = -113463.11
$ws.Range("H79").Value = 113778.11
$ws.Range("I79").Value = 113778.11
$ws.Range("K79").Value = 113778.11
$ws.Range("M79").Value = -112686.11
$ws.Range("H112").Value = 1457.0952
$ws.Range("J112").Value = 1515.7368
$ws.Range("L112").Value = 4547.2104
$ws.Range("N112").Value = -6763.2104
$ws.Range("H116").Value = 7527.5
$ws.Range("I116").Value = 1052
$ws.Range("J116").Value = 14003
$ws.Range("K116").Value = 1052
$ws.Range("L116").Value = 14003
$ws.Range("M116").Value = 2390
$ws.Range("N116").Value = -20887
$ws.Range("H121").Value = 785.7742
$ws.Range("J121").Value = 812.8214
$ws.Range("L121").Value = 2438.4642
$ws.Range("N121").Value = -5932.4642
$ws.Range("H129").Value = 78706.08
$ws.Range("I129").Value = 476
$ws.Range("K129").Value = 1428
$ws.Range("M129").Value = 3572

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 27778592
$ws.Range("I45").Value = 37037724
$ws.Range("J45").Value = 1200
$ws.Range("K45").Value = 37037724
$ws.Range("L45").Value = 1200
$ws.Range("M45").Value = -37037347
$ws.Range("N45").Value = -1954
$ws.Range("H61").Value = 2595.8667
$ws.Range("I61").Value = 2217.111
$ws.Range("K61").Value = 2217.111
$ws.Range("M61").Value = -2005.111
$ws.Range("H132").Value = 5625.224
$ws.Range("I132").Value = 6844.564
$ws.Range("J132").Value = 3122.3684
$ws.Range("K132").Value = 20533.692
$ws.Range("L132").Value = 9367.1052
$ws.Range("M132").Value = -18003.692
$ws.Range("N132").Value = -14427.1052
$ws.Range("H136").Value = 2595.8667
$ws.Range("I136").Value = 2217.111
$ws.Range("K136").Value = 6651.333
$ws.Range("M136").Value = -4101.333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2384.28
$ws.Range("I86").Value = 2250.3572
$ws.Range("J86").Value = 2554.7273
$ws.Range("K86").Value = 2250.3572
$ws.Range("L86").Value = 2554.7273
$ws.Range("M86").Value = -1127.3572
$ws.Range("N86").Value = -4800.7273
$ws.Range("H89").Value = 2384.28
$ws.Range("I89").Value = 2250.3572
$ws.Range("J89").Value = 2554.7273
$ws.Range("K89").Value = 11251.786
$ws.Range("L89").Value = 12773.6365
$ws.Range("M89").Value = -5635.786
$ws.Range("N89").Value = -24005.6365
$ws.Range("H99").Value = 688.1667
$ws.Range("I99").Value = 688.1667
$ws.Range("K99").Value = 688.1667
$ws.Range("M99").Value = 809.8333
$ws.Range("H107").Value = 1581.9
$ws.Range("I107").Value = 1272.125
$ws.Range("J107").Value = 2821
$ws.Range("K107").Value = 1272.125
$ws.Range("L107").Value = 2821
$ws.Range("M107").Value = 647.875
$ws.Range("N107").Value = -6661
$ws.Range("H134").Value = 5305.1143
$ws.Range("I134").Value = 6784.864
$ws.Range("J134").Value = 2800.923
$ws.Range("K134").Value = 20354.592
$ws.Range("L134").Value = 8402.769
$ws.Range("M134").Value = -17819.592
$ws.Range("N134").Value = -13472.769

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 142859150
$ws.Range("I86").Value = 166668260
$ws.Range("J86").Value = 4500
$ws.Range("K86").Value = 166668260
$ws.Range("L86").Value = 4500
$ws.Range("M86").Value = -166667137
$ws.Range("N86").Value = -6746
$ws.Range("H89").Value = 142859150
$ws.Range("I89").Value = 166668260
$ws.Range("J89").Value = 4500
$ws.Range("K89").Value = 833341300
$ws.Range("L89").Value = 22500
$ws.Range("M89").Value = -833335684
$ws.Range("N89").Value = -33732

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()
$ws.Range("H116").Value = 12199.8
$ws.Range("I116").Value = 18499.834
$ws.Range("K116").Value = 55499.50199999999
$ws.Range("M116").Value = -52057.50199999999
$ws.Range("H139").Value = 33334794
$ws.Range("I139").Value = 35715550
$ws.Range("J139").Value = 4200
$ws.Range("K139").Value = 107146650
$ws.Range("L139").Value = 12600
$ws.Range("M139").Value = -107141510
$ws.Range("N139").Value = -22880
$ws.Range("H140").Value = 2376.8975
$ws.Range("J140").Value = 2749.6155
$ws.Range("L140").Value = 8248.8465
$ws.Range("N140").Value = -18608.8465
$ws.Range("H141").Value = 3174.92
$ws.Range("I141").Value = 2877.2632
$ws.Range("J141").Value = 4117.5
$ws.Range("K141").Value = 8631.7896
$ws.Range("L141").Value = 12352.5
$ws.Range("M141").Value = -3451.7896
$ws.Range("N141").Value = -22712.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2188.3333
$ws.Range("I97").Value = 2090.8333
$ws.Range("J97").Value = 2383.3333
$ws.Range("K97").Value = 2090.8333
$ws.Range("L97").Value = 2383.3333
$ws.Range("M97").Value = -1594.8333
$ws.Range("N97").Value = -3375.3333
$ws.Range("H126").Value = 3072.4546
$ws.Range("I126").Value = 3157.5
$ws.Range("K126").Value = 9472.5
$ws.Range("M126").Value = -7002.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2270.4644
$ws.Range("I40").Value = 2304.9167
$ws.Range("J40").Value = 2063.75
$ws.Range("K40").Value = 2304.9167
$ws.Range("L40").Value = 2063.75
$ws.Range("M40").Value = -2168.9167
$ws.Range("N40").Value = -2335.75
$ws.Range("H46").Value = 1247.0358
$ws.Range("I46").Value = 697
$ws.Range("J46").Value = 1467.05
$ws.Range("K46").Value = 697
$ws.Range("L46").Value = 1467.05
$ws.Range("M46").Value = -509
$ws.Range("N46").Value = -1843.05
$ws.Range("H61").Value = 15152652
$ws.Range("I61").Value = 1210
$ws.Range("J61").Value = 33334382
$ws.Range("K61").Value = 1210
$ws.Range("L61").Value = 33334382
$ws.Range("M61").Value = -1008
$ws.Range("N61").Value = -33334786
$ws.Range("H93").Value = 1335.7273
$ws.Range("I93").Value = 1284.7142
$ws.Range("J93").Value = 1425
$ws.Range("K93").Value = 1284.7142
$ws.Range("L93").Value = 1425
$ws.Range("M93").Value = -36.71419999999989
$ws.Range("N93").Value = -3921
$ws.Range("H113").Value = 15152652
$ws.Range("I113").Value = 1210
$ws.Range("J113").Value = 33334382
$ws.Range("K113").Value = 1210
$ws.Range("L113").Value = 33334382
$ws.Range("M113").Value = 960
$ws.Range("N113").Value = -33338722
$ws.Range("H122").Value = 5084.5386
$ws.Range("I122").Value = 9240.799999999999
$ws.Range("J122").Value = 2486.875
$ws.Range("K122").Value = 27722.4
$ws.Range("L122").Value = 7460.625
$ws.Range("M122").Value = -25272.4
$ws.Range("N122").Value = -12360.625
$ws.Range("H132").Value = 9760.73
$ws.Range("I132").Value = 15570.286
$ws.Range("J132").Value = 2982.9167
$ws.Range("K132").Value = 46710.858
$ws.Range("L132").Value = 8948.750100000001
$ws.Range("M132").Value = -44180.858
$ws.Range("N132").Value = -14008.7501
